# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: insert a new day column (09-nov) before the
#    01-oct. column (column DN), shifting all following day columns
#    one to the right (DN..ER -> DO..ES). New column filled with the
#    same "-" placeholder used for the other not-yet-available future
#    days.
#  - "Gaz" and "CO2" sheets: append a new daily row (2025-11-07).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Prix Spot" - insert new day column at DN (shifts right)
# ---------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DN; existing DN..ER shift to DO..ES
$wsPrix.Columns("DN").Insert()

# Header (row 1) for the newly inserted column
$wsPrix.Range("DN1").Value = "09-nov"

# Data rows 2-25: no data yet for this day -> placeholder dash,
# matching the other not-yet-available day columns (DF:DM)
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 118).Value = "-"
}

# ---------------------------------------------------------------
# Sheet "Gaz" - append new daily row
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text (not date) parsing for the date-looking label, then drop the
# temporary number format so the cell keeps the workbook's default style.
$wsGaz.Range("A146").NumberFormat = "@"
$wsGaz.Range("A146").Value = "2025-11-07"
$wsGaz.Range("A146").ClearFormats()
$wsGaz.Range("B146").Value = 29.74

# ---------------------------------------------------------------
# Sheet "CO2" - append new daily row
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A146").NumberFormat = "@"
$wsCo2.Range("A146").Value = "2025-11-07"
$wsCo2.Range("A146").ClearFormats()
$wsCo2.Range("B146").Value = 79.36
